$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "2018 LEAVE CREDITS": monthly PERIOD dates in column A (rows 79-105)
# move from the 1st-of-month to the last day of the month, and the newly
# posted months (rows 84 & 85, i.e. Jun/Jul-2023) get an EARNED value.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("2018 LEAVE CREDITS")

$periods1 = @{
    79  = @(2023,1,31)
    80  = @(2023,2,28)
    81  = @(2023,3,31)
    82  = @(2023,4,30)
    83  = @(2023,5,31)
    84  = @(2023,6,30)
    85  = @(2023,7,31)
    86  = @(2023,8,31)
    87  = @(2023,9,30)
    88  = @(2023,10,31)
    89  = @(2023,11,30)
    90  = @(2023,12,31)
    91  = @(2024,1,31)
    92  = @(2024,2,29)
    93  = @(2024,3,31)
    94  = @(2024,4,30)
    95  = @(2024,5,31)
    96  = @(2024,6,30)
    97  = @(2024,7,31)
    98  = @(2024,8,31)
    99  = @(2024,9,30)
    100 = @(2024,10,31)
    101 = @(2024,11,30)
    102 = @(2024,12,31)
    103 = @(2025,1,31)
    104 = @(2025,2,28)
    105 = @(2025,3,31)
}

foreach ($row in $periods1.Keys) {
    $ymd = $periods1[$row]
    $d = Get-Date -Year $ymd[0] -Month $ymd[1] -Day $ymd[2] -Hour 0 -Minute 0 -Second 0
    $ws1.Range("A$row").Value = $d
}

# New leave credits earned for Jun/Jul-2023 rows (1.25 each)
$ws1.Range("C84").Value = 1.25
$ws1.Range("C85").Value = 1.25

# ---------------------------------------------------------------------------
# Sheet "2017 LEAVE BALANCE": new leave entries, rows 29-31
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# Row 29: SL(1-0-0), 1 day absence with pay, covering 7/1 - 7/14/2023
$ws2.Range("A29").Value = (Get-Date -Year 2023 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws2.Range("B29").Value = "SL(1-0-0)"
$ws2.Range("H29").Value = 1
$ws2.Range("K29").Value = (Get-Date -Year 2023 -Month 7 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws2.Range("K29").NumberFormat = "m/d/yyyy"

# Row 30: VL(3-0-0), 3 days absence with pay, covering 7/21,25,28/2023
$ws2.Range("B30").Value = "VL(3-0-0)"
$ws2.Range("D30").Value = 3
$ws2.Range("K30").Value = "7/21,25,28/2023"

# Row 31: VL(10-0-0), 10 days absence with pay, covering 8/1,4,8,11,15,18,22,25,29,31
$ws2.Range("A31").Value = (Get-Date -Year 2023 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0)
$ws2.Range("B31").Value = "VL(10-0-0)"
$ws2.Range("D31").Value = 10
$ws2.Range("K31").Value = "8/1,4,8,11,15,18,22,25,29,31"

# ---------------------------------------------------------------------------
# UI state: active sheet / selection, matching the saved view
# ---------------------------------------------------------------------------
$ws2.ListObjects.Item("Table1").DataBodyRange.Select()

$ws1.Activate()
$ws1.Range("C84:C85").Select()
